$d = $word.ActiveDocument

# 1. Small text correction: drop the "(yfinance)" aside from the free data sources line.
$d.Content.Find.Execute(
    "Free: Yahoo Finance (yfinance), Alpha Vantage, Quandl.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Free: Yahoo Finance, Alpha Vantage, Quandl.", 2)

# 2. Append a new "Sentiment Analysis APIs" bullet (top level, ilvl 0) plus four
#    sub-bullets (ilvl 1), each holding a link to one of the new APIs.

# --- top-level heading bullet ---
$lastPara = $d.Paragraphs.Last
$insertRange = $lastPara.Range
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Range.Text = "Sentiment Analysis APIs"
$headingPara.Range.ListFormat.ListLevelNumber = 1
$headingPara.LeftIndent = 36

$urls = @(
    "https://tradestie.com/apps/reddit/api/",
    "https://housestockwatcher.com/api",
    "https://stocknewsapi.com/?ref=apilist.fun",
    "https://www.alphavantage.co/?ref=apilist.fun"
)

$prevPara = $headingPara
foreach ($url in $urls) {
    $lineText = $url + " "

    $prevRange = $prevPara.Range
    $prevRange.Collapse(0)
    $prevRange.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Range.Text = $lineText
    $p.Range.ListFormat.ListLevelNumber = 2
    $p.LeftIndent = 72

    $pRange = $p.Range
    $urlStart = $pRange.Start
    $urlEnd = $urlStart + $url.Length
    $urlRange = $d.Range($urlStart, $urlEnd)

    $link = $d.Hyperlinks.Add($urlRange, $url, "", "", $url)
    $urlRange.Font.Color = 13391121
    $urlRange.Font.Underline = 1

    $prevPara = $p
}
